$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to append: (A value, B value) pairs, continuing the existing series.
$newData = @(
    @(204, 0.2366255144032922),
    @(205, 0.3099513655069209),
    @(206, 0.1411522633744856),
    @(207, 0.3518518518518517),
    @(208, 0.5349794238683127),
    @(209, 0.2342739564961787),
    @(210, 0.7798353909465019),
    @(211, 0.5855967078189299),
    @(212, 0.5102880658436213),
    @(213, 0.5164609053497942),
    @(214, 0.368312757201646),
    @(215, 0.5164609053497942)
)

# Existing data occupies rows 2..205 (header row 1). New rows start right after.
$startRow = 206

# Use the same formatting as the existing column A cells (e.g. A205) for the new A cells.
$templateA = $ws.Range("A205")
$templateA.Copy()

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $aVal = $newData[$i][0]
    $bVal = $newData[$i][1]

    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)

    $aCell.Value = $aVal
    $bCell.Value = $bVal

    # Match formatting (bold, centered, bordered) of existing A column cells.
    $aCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0
